$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, pushing the existing rows 14-46 down to 15-47.
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly entry. The
# "static" columns (A,B,C,E,F,G,H,I,N,O,Q,R) carry the same values as
# every other row in this sheet; only the date/volume/price columns
# (D,J,K,L,M,P) hold data specific to this record.
$ws.Range("A14").Value() = 3
$ws.Range("B14").Value() = "Femacal de La Calera"
$ws.Range("C14").Value() = "Coquimbo"
$ws.Range("D14").Value() = "2022-06-30"
$ws.Range("E14").Value() = 5
$ws.Range("F14").Value() = 100112035
$ws.Range("G14").Value() = "Bruselas (repollito)"
$ws.Range("H14").Value() = "Sin especificar"
$ws.Range("I14").Value() = "Primera"
$ws.Range("J14").Value() = 85
$ws.Range("K14").Value() = 15000
$ws.Range("L14").Value() = 16000
$ws.Range("M14").Value() = 15529
$ws.Range("N14").Value() = "$/malla 15 kilos"
$ws.Range("O14").Value() = "Provincia de Quillota"
$ws.Range("P14").Value() = 1035
$ws.Range("Q14").Value() = 15
$ws.Range("R14").Value() = "Hortaliza"
